# Apply updated crypto price/volume data to Sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.723.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.002.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.523"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.999.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.492.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.621.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.002.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.22%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0833"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "397.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0355"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.269"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.723.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
